$d = $word.ActiveDocument

# Locate the paragraph that currently holds the lone "_GoBack" bookmark.
# It is the first empty paragraph immediately after the paragraph containing
# "Neste sentido, a interface do orcamentista" text.
$target = $null
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Neste sentido*") {
        $target = $d.Paragraphs.Item($i + 1)
        break
    }
}

if ($null -eq $target) {
    throw "Could not locate the anchor paragraph (bookmark paragraph) to replace."
}

$xml = '    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
    </w:p>
    <w:tbl>
      <w:tblPr>
        <w:tblStyle w:val="Tabelacomgrade"/>
        <w:tblW w:w="0" w:type="auto"/>
        <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="8720"/>
      </w:tblGrid>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="8720" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t xml:space="preserve">1 - </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t>Com relação à orçamentação de uma obra, analise as afirmativas a seguir e assinale V para a verdadeira e F para a falsa.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t>(    ) Os custos com a Administração Local são despesas indiretas e devem ser inseridas no cálculo do BDI da obra.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t>(    ) Os custos com a mobilização e desmobilização da obra fazem parte dos custos diretos da obra e devem constar na planilha orçamentária.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t>(    ) Os custos com a instalação do canteiro e acampamento da obra são custos indiretos e devem ser inseridos no cálculo do BDI da obra.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t>As afirmativas são, na ordem apresentada, respectivamente,</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t xml:space="preserve">R: </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t>F – V – F</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t xml:space="preserve">2 - </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t>Sobre o SINAPI, é correto afirmar que:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t>A Caixa Econômica Federal e o Instituto Brasileiro de Geografia e Estatística (IBGE) são responsáveis pela aferição e divulgação de preços de serviços e insumos para utilização em orçamentos.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t xml:space="preserve">3 - </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t>Todo orçamento é composto por custos diretos, custos indiretos e o BDI − Benefícios e Despesas Indiretas. Não corresponde a uma parcela de custo direto:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t xml:space="preserve">R: </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t>Administração central.</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t>.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t xml:space="preserve">4 - </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t>Em função do grau de suficiência das informações e sua correta compreensão, os produtos de orçamento podem apresentar margem de erro. Dentre as alternativas abaixo, aponte a incorreta:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t xml:space="preserve">R: </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t>Avaliação, margem de erro de 20 a 35%;</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t xml:space="preserve">5 - </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t>Sobre o processo de orçamentação direcionado ao Orçamento Operacional, leia as alternativas abaixo e sinalize a incorreta:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t xml:space="preserve">R: </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
              <w:t>Define Serviço como o conjunto de tarefas que resultam em uma parte funcional da obra.</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="6562"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
          <w:noProof/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:tbl>
      <w:tblPr>
        <w:tblStyle w:val="Tabelacomgrade"/>
        <w:tblW w:w="0" w:type="auto"/>
        <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="8720"/>
      </w:tblGrid>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="8720" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:tabs>
                <w:tab w:val="left" w:pos="6562"/>
              </w:tabs>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
                <w:noProof/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
              </w:rPr>
            </w:pPr>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="6562"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
          <w:noProof/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="6562"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
          <w:noProof/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="6562"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
          <w:noProof/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="6562"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
          <w:noProof/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="6562"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
          <w:noProof/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="6562"/>
        </w:tabs>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>
          <w:noProof/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:lang w:val="pt-BR" w:bidi="ar-SA"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'

$target.Range.InsertXML($xml)

Write-Output "Replacement complete"
